$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the date separator from "/" to "-" for all date cells in column A (rows 3-21).
# Some of these strings look like valid dates to Excel's auto-detection (e.g. "01-08-2022"),
# which would otherwise silently convert them to a date serial number. Forcing the cell to
# Text format before assigning (and resetting the style afterwards) keeps them as plain text,
# matching the original inline-string representation.
for ($r = 3; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $newText = $cell.Text -replace '/', '-'
    $cell.NumberFormat = "@"
    $cell.Value = $newText
    $cell.Style = "Normal"
}

# Update attendance counts for row 3 (28-07-2022)
$ws.Cells.Item(3, 4).Value = 1   # D3: 0 -> 1
$ws.Cells.Item(3, 7).Value = 1   # G3: 0 -> 1

# Update attendance counts for row 10 (22-08-2022)
$ws.Cells.Item(10, 4).Value = 1  # D10: 0 -> 1
$ws.Cells.Item(10, 5).Value = 1  # E10: 0 -> 1
$ws.Cells.Item(10, 8).Value = 0  # H10: 1 -> 0
